$wb = $excel.ActiveWorkbook

$wsAbout = $wb.Worksheets.Item("About")
$wsData = $wb.Worksheets.Item("Boundaries and methane sources")

$newVersion = "mines - version 1.0.0 (Feb 3 2026) (built on February 03 2026 10.14.00 EST)"

# Row 2: Version line
$wsAbout.Range("A2").Value = "Version: $newVersion"

# Row 6: Recommended Citation line
$wsAbout.Range("A6").Value = "Recommended Citation:  `"Global Energy Monitor, Coal mine boundaries and methane sources for Shanxi Jinyuan Coal Mine, China, M0330, version '$newVersion'. (See the CC license for attribution requirements if sharing or adapting the data set.)"

# Column S (build_version) rows 2-14 on the data sheet
for ($r = 2; $r -le 14; $r++) {
    $wsData.Cells.Item($r, 19).Value = $newVersion
}
